$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.206.95'
$ws.Range("D2").Style = $origStyle
$origStyle = $ws.Range("E2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.83%  '
$ws.Range("E2").Style = $origStyle

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.796.33'
$ws.Range("D3").Style = $origStyle
$origStyle = $ws.Range("E3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.88%  '
$ws.Range("E3").Style = $origStyle

$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = $origStyle
$origStyle = $ws.Range("E4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E4").Style = $origStyle

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '339.01'
$ws.Range("D5").Style = $origStyle
$origStyle = $ws.Range("E5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("E5").Style = $origStyle

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9990'
$ws.Range("D6").Style = $origStyle
$origStyle = $ws.Range("E6").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E6").Style = $origStyle

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4565'
$ws.Range("D7").Style = $origStyle
$origStyle = $ws.Range("E7").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +21.07%  '
$ws.Range("E7").Style = $origStyle

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3602'
$ws.Range("D8").Style = $origStyle
$origStyle = $ws.Range("E8").Style
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +6.69%  '
$ws.Range("E8").Style = $origStyle

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.46'
$ws.Range("D9").Style = $origStyle
$origStyle = $ws.Range("E9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.38%  '
$ws.Range("E9").Style = $origStyle

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.139'
$ws.Range("D10").Style = $origStyle
$origStyle = $ws.Range("E10").Style
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.55%  '
$ws.Range("E10").Style = $origStyle

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07498'
$ws.Range("D11").Style = $origStyle
$origStyle = $ws.Range("E11").Style
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.82%  '
$ws.Range("E11").Style = $origStyle

$origStyle = $ws.Range("E12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("E12").Style = $origStyle

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.35'
$ws.Range("D13").Style = $origStyle
$origStyle = $ws.Range("E13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.09%  '
$ws.Range("E13").Style = $origStyle

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.210'
$ws.Range("D14").Style = $origStyle
$origStyle = $ws.Range("E14").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.04%  '
$ws.Range("E14").Style = $origStyle

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.236'
$ws.Range("D15").Style = $origStyle
$origStyle = $ws.Range("E15").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.34%  '
$ws.Range("E15").Style = $origStyle

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.793.12'
$ws.Range("D16").Style = $origStyle
$origStyle = $ws.Range("E16").Style
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.65%  '
$ws.Range("E16").Style = $origStyle

$origStyle = $ws.Range("E17").Style
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.53%  '
$ws.Range("E17").Style = $origStyle

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06691'
$ws.Range("D18").Style = $origStyle
$origStyle = $ws.Range("E18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.77%  '
$ws.Range("E18").Style = $origStyle

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.16'
$ws.Range("D19").Style = $origStyle
$origStyle = $ws.Range("E19").Style
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.48%  '
$ws.Range("E19").Style = $origStyle

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9995'
$ws.Range("D20").Style = $origStyle
$origStyle = $ws.Range("E20").Style
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("E20").Style = $origStyle

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.20'
$ws.Range("D21").Style = $origStyle
$origStyle = $ws.Range("E21").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.38%  '
$ws.Range("E21").Style = $origStyle

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.368'
$ws.Range("D22").Style = $origStyle
$origStyle = $ws.Range("E22").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.51%  '
$ws.Range("E22").Style = $origStyle

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.240.13'
$ws.Range("D23").Style = $origStyle
$origStyle = $ws.Range("E23").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.92%  '
$ws.Range("E23").Style = $origStyle

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.87'
$ws.Range("D24").Style = $origStyle
$origStyle = $ws.Range("E24").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.38%  '
$ws.Range("E24").Style = $origStyle

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.370'
$ws.Range("D25").Style = $origStyle
$origStyle = $ws.Range("E25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.02%  '
$ws.Range("E25").Style = $origStyle

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.37'
$ws.Range("D26").Style = $origStyle
$origStyle = $ws.Range("E26").Style
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.89%  '
$ws.Range("E26").Style = $origStyle

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.63'
$ws.Range("D27").Style = $origStyle
$origStyle = $ws.Range("E27").Style
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.29%  '
$ws.Range("E27").Style = $origStyle

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.384'
$ws.Range("D28").Style = $origStyle
$origStyle = $ws.Range("E28").Style
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.27%  '
$ws.Range("E28").Style = $origStyle

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.994.32'
$ws.Range("D29").Style = $origStyle
$origStyle = $ws.Range("E29").Style
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.56%  '
$ws.Range("E29").Style = $origStyle

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.264'
$ws.Range("D30").Style = $origStyle
$origStyle = $ws.Range("E30").Style
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("E30").Style = $origStyle

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '132.10'
$ws.Range("D31").Style = $origStyle
$origStyle = $ws.Range("E31").Style
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.60%  '
$ws.Range("E31").Style = $origStyle

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.071'
$ws.Range("D32").Style = $origStyle
$origStyle = $ws.Range("E32").Style
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.37%  '
$ws.Range("E32").Style = $origStyle

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.865'
$ws.Range("D33").Style = $origStyle
$origStyle = $ws.Range("E33").Style
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.48%  '
$ws.Range("E33").Style = $origStyle

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09440'
$ws.Range("D34").Style = $origStyle
$origStyle = $ws.Range("E34").Style
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +7.28%  '
$ws.Range("E34").Style = $origStyle

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02363'
$ws.Range("D35").Style = $origStyle
$origStyle = $ws.Range("E35").Style
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.67%  '
$ws.Range("E35").Style = $origStyle

$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.05'
$ws.Range("D36").Style = $origStyle
$origStyle = $ws.Range("E36").Style
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.66%  '
$ws.Range("E36").Style = $origStyle

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06262'
$ws.Range("D37").Style = $origStyle
$origStyle = $ws.Range("E37").Style
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.63%  '
$ws.Range("E37").Style = $origStyle

$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6619'
$ws.Range("D38").Style = $origStyle
$origStyle = $ws.Range("E38").Style
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.34%  '
$ws.Range("E38").Style = $origStyle

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.159'
$ws.Range("D39").Style = $origStyle
$origStyle = $ws.Range("E39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("E39").Style = $origStyle

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2154'
$ws.Range("D40").Style = $origStyle
$origStyle = $ws.Range("E40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.92%  '
$ws.Range("E40").Style = $origStyle

$origStyle = $ws.Range("E41").Style
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.29%  '
$ws.Range("E41").Style = $origStyle

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.214'
$ws.Range("D42").Style = $origStyle
$origStyle = $ws.Range("E42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("E42").Style = $origStyle

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.043'
$ws.Range("D43").Style = $origStyle
$origStyle = $ws.Range("E43").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.24%  '
$ws.Range("E43").Style = $origStyle

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9994'
$ws.Range("D44").Style = $origStyle
$origStyle = $ws.Range("E44").Style
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("E44").Style = $origStyle

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.97'
$ws.Range("D45").Style = $origStyle
$origStyle = $ws.Range("E45").Style
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.54%  '
$ws.Range("E45").Style = $origStyle

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.869'
$ws.Range("D46").Style = $origStyle

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6061'
$ws.Range("D47").Style = $origStyle
$origStyle = $ws.Range("E47").Style
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("E47").Style = $origStyle

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.17'
$ws.Range("D48").Style = $origStyle
$origStyle = $ws.Range("E48").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("E48").Style = $origStyle

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.018'
$ws.Range("D49").Style = $origStyle
$origStyle = $ws.Range("E49").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.27%  '
$ws.Range("E49").Style = $origStyle

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07090'
$ws.Range("D50").Style = $origStyle
$origStyle = $ws.Range("E50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.04%  '
$ws.Range("E50").Style = $origStyle

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.161'
$ws.Range("D51").Style = $origStyle
$origStyle = $ws.Range("E51").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.80%  '
$ws.Range("E51").Style = $origStyle
